# Aggiornamento 15, 16, 17 marzo: append three new daily rows (227-229)
# after the last existing data row (226), extending the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date column's formatting (style incl. date number format, border,
# font, alignment) from the last existing row down onto the new rows so the
# new A-column cells keep the same look as the rest of the date column.
$ws.Range("A226").Copy() | Out-Null
$ws.Range("A227:A229").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 1
$ws.Range("C227").Value = 7
$ws.Range("D227").Value = 106.6098081023454

$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 2
$ws.Range("C228").Value = 8
$ws.Range("D228").Value = 121.8397806883948

$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 1
$ws.Range("C229").Value = 5
$ws.Range("D229").Value = 76.14986293024673
